# Automatic update of files.
#
# Applies the diff to rows 3, 4 and 7 of the active sheet:
#   - Id (A) changed to new observation ids
#   - Noggrannhet (S) 25 -> 10
#   - Startdatum/Starttid/Slutdatum/Sluttid (Y/Z/AA/AB) updated
#   - Substrat-beskrivning (AO7) "salg" -> "Salg" (capitalised)
#   - Rapportor/Observatorer/Projektnamn (AW/AX/AY) updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date/time columns (Y, Z, AA, AB) hold plain text values (e.g. "2022-08-14",
# "13:12") rather than real Excel dates/times. Assigning a date-looking string
# straight to .Value makes Excel auto-convert it into a date serial, so each
# cell is first forced to Text format, written, and then its style is reset
# back to Normal (keeping the stored value textual without leaving a lasting
# custom number format on the cell).
$dateTimeCells = @("Y3","Z3","AA3","AB3","Y4","Z4","AA4","AB4","Y7","Z7","AA7","AB7")
foreach ($addr in $dateTimeCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 3
$ws.Range("A3").Value = 104025332
$ws.Range("S3").Value = 10
$ws.Range("Y3").Value = "2022-08-14"
$ws.Range("Z3").Value = "13:12"
$ws.Range("AA3").Value = "2022-08-14"
$ws.Range("AB3").Value = "13:12"
$ws.Range("AW3").Value = "Johan Staaf"
$ws.Range("AX3").Value = "Via Johan Staaf"
$ws.Range("AY3").Value = "LstZ inventering av skogliga värdetrakter 2022"

# Row 4
$ws.Range("A4").Value = 104025311
$ws.Range("S4").Value = 10
$ws.Range("Y4").Value = "2022-08-14"
$ws.Range("Z4").Value = "13:12"
$ws.Range("AA4").Value = "2022-08-14"
$ws.Range("AB4").Value = "13:12"
$ws.Range("AW4").Value = "Johan Staaf"
$ws.Range("AX4").Value = "Via Johan Staaf"
$ws.Range("AY4").Value = "LstZ inventering av skogliga värdetrakter 2022"

# Row 7
$ws.Range("A7").Value = 104025346
$ws.Range("S7").Value = 10
$ws.Range("Y7").Value = "2022-08-14"
$ws.Range("Z7").Value = "13:31"
$ws.Range("AA7").Value = "2022-08-14"
$ws.Range("AB7").Value = "13:31"
$ws.Range("AO7").Value = "Sälg"
$ws.Range("AW7").Value = "Johan Staaf"
$ws.Range("AX7").Value = "Via Johan Staaf"
$ws.Range("AY7").Value = "LstZ inventering av skogliga värdetrakter 2022"

# Restore the default (General) style on the date/time cells so the only
# observable change is the cell's text content, matching the source diff.
$ws.Range("Y3:AB4").Style = "Normal"
$ws.Range("Y7:AB7").Style = "Normal"
